$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '67.340.99'
$ws.Range("E2").Value = '  +5.09%  '

# Row 3
$ws.Range("D3").Value = '3.249.18'
$ws.Range("E3").Value = '  +2.71%  '

# Row 4
$ws.Range("E4").Value = '  +0.02%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '576.54'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.43%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '179.35'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +6.70%  '

# Row 7
$ws.Range("E7").Value = '  +0.01%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.600'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.81%  '

# Row 9
$ws.Range("D9").Value = '3.247.64'
$ws.Range("E9").Value = '  +2.76%  '

# Row 10
$ws.Range("E10").Value = '  +4.45%  '

# Row 11
$ws.Range("E11").Value = '  +2.06%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.414'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +5.02%  '

# Row 13
$ws.Range("D13").Value = '3.810.67'
$ws.Range("E13").Value = '  +2.75%  '

# Row 14
$ws.Range("E14").Value = '  +0.63%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '27.97'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.65%  '

# Row 16
$ws.Range("D16").Value = '67.257.45'

# Row 17
$ws.Range("E17").Value = '  +3.16%  '

# Row 18
$ws.Range("D18").Value = '3.252.16'
$ws.Range("E18").Value = '  +2.91%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.87'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.57%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.34'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +3.17%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '375.70'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +6.95%  '

# Row 22
$ws.Range("E22").Value = '  +6.26%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.00'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.10%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '71.14'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +3.78%  '

# Row 25
$ws.Range("E25").Value = '  +1.72%  '

# Row 26
$ws.Range("D26").Value = '3.388.19'
$ws.Range("E26").Value = '  +2.61%  '

# Row 27
$ws.Range("E27").Value = '  -0.57%  '

# Row 28
$ws.Range("E28").Value = '  +4.77%  '

# Row 29
$ws.Range("E29").Value = '  +1.99%  '

# Row 30
$ws.Range("E30").Value = '  +0.26%  '

# Row 31
$ws.Range("E31").Value = '  +4.66%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.62'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.18%  '

# Row 33
$ws.Range("E33").Value = '  +2.74%  '

# Row 34
$ws.Range("E34").Value = '  +0.06%  '

# Row 35
$ws.Range("E35").Value = '  +5.99%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.83'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +3.20%  '

# Row 37
$ws.Range("B37").Value = 'ImmutableX'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.49'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +4.01%  '

# Row 38
$ws.Range("B38").Value = 'Monero'
$ws.Range("C38").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '161.79'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +5.39%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.857'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +5.57%  '

# Row 40
$ws.Range("E40").Value = '  +10.08%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.86'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +14.93%  '

# Row 42
$ws.Range("E42").Value = '  +4.68%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.63'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +5.47%  '

# Row 44
$ws.Range("D44").Value = '2.758.98'
$ws.Range("E44").Value = '  +6.32%  '

# Row 45
$ws.Range("E45").Value = '  +5.61%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '352.04'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +10.92%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '25.78'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +9.08%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '40.46'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.83%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0672'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.43%  '

# Row 50
$ws.Range("E50").Value = '  +4.32%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.102'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.59%  '
